$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.702.07"
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").Value = "2.294.37"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'303.99"
$ws.Range("E5").Value = "  +1.28%  "

$ws.Range("D6").Value = "'96.29"
$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("D7").Value = "'0.506"
$ws.Range("E7").Value = "  -2.18%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -2.34%  "

$ws.Range("D10").Value = "'34.85"
$ws.Range("E10").Value = "  -3.24%  "

$ws.Range("E11").Value = "  -0.88%  "

$ws.Range("D12").Value = "'18.63"
$ws.Range("E12").Value = "  +4.63%  "

$ws.Range("E13").Value = "  +1.80%  "

$ws.Range("D14").Value = "'6.85"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("D15").Value = "2.649.69"
$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").Value = "2.284.25"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").Value = "42.620.94"
$ws.Range("E18").Value = "  -0.78%  "

$ws.Range("D19").Value = "'12.98"
$ws.Range("E19").Value = "  +2.77%  "

$ws.Range("E20").Value = "  -1.59%  "

$ws.Range("D21").Value = "'5.98"
$ws.Range("E21").Value = "  -2.07%  "

$ws.Range("D22").Value = "'67.18"
$ws.Range("E22").Value = "  -1.51%  "

$ws.Range("D23").Value = "'235.90"
$ws.Range("E23").Value = "  -2.54%  "

$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("D26").Value = "'2.39"
$ws.Range("E26").Value = "  -1.64%  "

$ws.Range("D27").Value = "'24.77"
$ws.Range("E27").Value = "  -1.74%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'167.18"
$ws.Range("E28").Value = "  +0.25%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.06"
$ws.Range("E29").Value = "  +0.72%  "

$ws.Range("D30").Value = "'8.99"
$ws.Range("E30").Value = "  -0.93%  "

$ws.Range("D31").Value = "'32.92"
$ws.Range("E31").Value = "  -0.45%  "

$ws.Range("D33").Value = "'17.78"
$ws.Range("E33").Value = "  +0.51%  "

$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").Value = "'4.45"
$ws.Range("E35").Value = "  -5.92%  "

$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("E39").Value = "  -1.47%  "

$ws.Range("E40").Value = "  -1.25%  "

$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = "  -3.01%  "

$ws.Range("D42").Value = "1.993.64"
$ws.Range("E42").Value = "  -0.57%  "

$ws.Range("D43").Value = "'0.0278"
$ws.Range("E43").Value = "  -3.04%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'10.25"
$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'18.26"
$ws.Range("E45").Value = "  +4.56%  "

$ws.Range("D46").Value = "'2.11"
$ws.Range("E46").Value = "  -2.41%  "

$ws.Range("D47").Value = "'2.76"
$ws.Range("E47").Value = "  -0.73%  "

$ws.Range("D48").Value = "'2.87"
$ws.Range("E48").Value = "  -1.69%  "

$ws.Range("D49").Value = "'53.66"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").Value = "2.513.89"
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'70.58"
$ws.Range("E51").Value = "  -2.98%  "
